# Generate Report for Archive
#
# 1. Flip the localization status from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2 for the
#    zh-cn / de-de columns, and the per-language "Status" column on
#    the zh-cn and de-de report sheets, zh-cn!C2 / de-de!C2).
# 2. Narrow the "Status"-related columns that used to be sized for the
#    longer "Ready for handoff" text down to fit "In Translation"
#    (Overview columns E:F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 13.4101845877511

# --- Overview sheet: zh-cn (E) / de-de (F) status cells + widths ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) cell + width ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) cell + width ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
